$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 10, pushing the
# existing rows 10-13 down to 11-14.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 45040
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 100112041
$ws.Cells.Item(10, 7).Value = "Fruto del paraíso"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 18000
$ws.Cells.Item(10, 14).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 1000
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
